$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.06446533333333333
$ws.Range("H2").Value = 0.193396
$ws.Range("I2").Value = 0.02693738696927793
$ws.Range("J2").Value = 0.02693738696927793
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.030956000000001
$ws.Range("N2").Value = 18.092868
$ws.Range("O2").Value = 0.364814105361131
$ws.Range("P2").Value = 0.3648141053611309
$ws.Range("Q2").Value = 0.3887875888586667
$ws.Range("R2").Value = 3.499088299728001
$ws.Range("S2").Value = 0.009827138727963715
$ws.Range("T2").Value = 0.009827138727963713
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.06446533333333333
$ws.Range("H3").Value = 0.193396
$ws.Range("I3").Value = 0.02693738696927793
$ws.Range("J3").Value = 0.02693738696927793
$ws.Range("O3").Value = 0.4107214552505144
$ws.Range("P3").Value = 0.4107214552505143
$ws.Range("Q3").Value = 0.4377117055857778
$ws.Range("R3").Value = 3.939405350272
$ws.Range("S3").Value = 0.01106376277666807
$ws.Range("T3").Value = 0.01106376277666807
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.06446533333333333
$ws.Range("H4").Value = 0.193396
$ws.Range("I4").Value = 0.02693738696927793
$ws.Range("J4").Value = 0.02693738696927793
$ws.Range("M4").Value = 3.710753333333333
$ws.Range("N4").Value = 11.13226
$ws.Range("O4").Value = 0.2244644393883547
$ws.Range("P4").Value = 0.2244644393883547
$ws.Range("Q4").Value = 0.2392149505511111
$ws.Range("R4").Value = 2.15293455496
$ws.Range("S4").Value = 0.006046485464646142
$ws.Range("T4").Value = 0.006046485464646141
$ws.Range("I5").Value = 0.7704314695358874
$ws.Range("J5").Value = 0.7704314695358874
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.030956000000001
$ws.Range("N5").Value = 18.092868
$ws.Range("O5").Value = 0.364814105361131
$ws.Range("P5").Value = 0.3648141053611309
$ws.Range("Q5").Value = 11.11964548615333
$ws.Range("R5").Value = 100.07680937538
$ws.Range("S5").Value = 0.2810642673007961
$ws.Range("T5").Value = 0.2810642673007961
$ws.Range("I6").Value = 0.7704314695358874
$ws.Range("J6").Value = 0.7704314695358874
$ws.Range("O6").Value = 0.4107214552505144
$ws.Range("P6").Value = 0.4107214552505143
$ws.Range("S6").Value = 0.316432734338572
$ws.Range("T6").Value = 0.3164327343385719
$ws.Range("I7").Value = 0.7704314695358874
$ws.Range("J7").Value = 0.7704314695358874
$ws.Range("M7").Value = 3.710753333333333
$ws.Range("N7").Value = 11.13226
$ws.Range("O7").Value = 0.2244644393883547
$ws.Range("P7").Value = 0.2244644393883547
$ws.Range("Q7").Value = 6.841744750455555
$ws.Range("R7").Value = 61.5757027541
$ws.Range("S7").Value = 0.1729344678965193
$ws.Range("T7").Value = 0.1729344678965192
$ws.Range("G8").Value = 0.4849276666666666
$ws.Range("H8").Value = 1.454783
$ws.Range("I8").Value = 0.2026311434948347
$ws.Range("J8").Value = 0.2026311434948347
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.030956000000001
$ws.Range("N8").Value = 18.092868
$ws.Range("O8").Value = 0.364814105361131
$ws.Range("P8").Value = 0.3648141053611309
$ws.Range("Q8").Value = 2.924577420849333
$ws.Range("R8").Value = 26.321196787644
$ws.Range("S8").Value = 0.07392269933237107
$ws.Range("T8").Value = 0.07392269933237106
$ws.Range("G9").Value = 0.4849276666666666
$ws.Range("H9").Value = 1.454783
$ws.Range("I9").Value = 0.2026311434948347
$ws.Range("J9").Value = 0.2026311434948347
$ws.Range("O9").Value = 0.4107214552505144
$ws.Range("P9").Value = 0.4107214552505143
$ws.Range("Q9").Value = 3.292599372206222
$ws.Range("R9").Value = 29.633394349856
$ws.Range("S9").Value = 0.08322495813527431
$ws.Range("T9").Value = 0.08322495813527428
$ws.Range("G10").Value = 0.4849276666666666
$ws.Range("H10").Value = 1.454783
$ws.Range("I10").Value = 0.2026311434948347
$ws.Range("J10").Value = 0.2026311434948347
$ws.Range("M10").Value = 3.710753333333333
$ws.Range("N10").Value = 11.13226
$ws.Range("O10").Value = 0.2244644393883547
$ws.Range("P10").Value = 0.2244644393883547
$ws.Range("Q10").Value = 1.799446955508889
$ws.Range("R10").Value = 16.19502259958
$ws.Range("S10").Value = 0.04548348602718933
$ws.Range("T10").Value = 0.04548348602718932
